$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-24 04:46:27"
$wsZhCn.Range("H4").Value = "2016-03-24 04:47:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-24 04:46:32"
$wsDeDe.Range("H4").Value = "2016-03-24 04:47:08"
